# satya code for payment and reusable components
#
# Adds three new "Payment" module test rows to the Run_Manager sheet and
# flips the "Execute" flag on the existing makePaymentwithExistingCreditCard
# row (110) from Yes to No, while the newly-appended
# verifyBankPaymentFormFieldsValidation row (113) is the one left set to run
# (Yes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Run_Manager")

# --- Row 110: makePaymentwithExistingCreditCard no longer set to execute ---
$ws.Cells.Item(110, 4).Value = "No"

# --- Prep rows 111-113: clone row 110's cell formatting (styles only) so
#     the new rows pick up the same column styles (center alignment,
#     quote-prefixed numeric-as-text cells, etc.) without touching styles.xml ---
$ws.Range("A110:F110").Copy()
$ws.Range("A111:F113").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 111: verifyPaymentInformationPageObjects ---
$ws.Cells.Item(111, 1).Value = "Payment"
$ws.Cells.Item(111, 2).Value = "verifyPaymentInformationPageObjects"
$ws.Cells.Item(111, 3).Value = "to verify payment information page objects"
$ws.Cells.Item(111, 4).Value = "No"
$ws.Cells.Item(111, 5).Value = "'1"
$ws.Cells.Item(111, 6).Value = "'1"

# --- Row 112: verifyBankPaymentFormFields ---
$ws.Cells.Item(112, 1).Value = "Payment"
$ws.Cells.Item(112, 2).Value = "verifyBankPaymentFormFields"
$ws.Cells.Item(112, 3).Value = "verify bank payment form fields"
$ws.Cells.Item(112, 4).Value = "No"
$ws.Cells.Item(112, 5).Value = "'1"
$ws.Cells.Item(112, 6).Value = "'1"

# --- Row 113: verifyBankPaymentFormFieldsValidation (set to execute) ---
$ws.Cells.Item(113, 1).Value = "Payment"
$ws.Cells.Item(113, 2).Value = "verifyBankPaymentFormFieldsValidation"
$ws.Cells.Item(113, 3).Value = "validation on bank payment form"
$ws.Cells.Item(113, 4).Value = "Yes"
$ws.Cells.Item(113, 5).Value = "'1"
$ws.Cells.Item(113, 6).Value = "'1"

# --- Update sheet view / selection to mirror the author's cursor position ---
$ws.Range("C116").Select()
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
